$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets = $wb.Worksheets.Item("Assets")

# Add the two new IMAP asset rows to the Assets sheet.
$wsAssets.Cells.Item(8, 1).Value = "IMAP_Port"
$wsAssets.Cells.Item(8, 2).Value = "IMAP_Port"
$wsAssets.Cells.Item(8, 3).Value = "LazyFramework"
$wsAssets.Cells.Item(8, 4).Value = "The port of the IMAP server."

$wsAssets.Cells.Item(9, 1).Value = "IMAP_Server"
$wsAssets.Cells.Item(9, 2).Value = "IMAP_Server"
$wsAssets.Cells.Item(9, 3).Value = "LazyFramework"
$wsAssets.Cells.Item(9, 4).Value = "The URL of the IMAP server."

# Update the selection saved on the Settings sheet (no longer the active tab).
$wsSettings.Activate()
$wsSettings.Range("C15").Select()

# Assets becomes the active/selected tab, with its own updated selection.
$wsAssets.Activate()
$wsAssets.Range("E8").Select()
